$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New simulated-subject rows (s16..s20) appended to the manifest table.
# All values are text (inlineStr) in the source workbook, including the
# numeric-looking ones, so we force text formatting before writing, then
# reset the cell style back to Normal so no extra number-format style
# lingers on the saved cells (matches the target which carries no `s=`
# style override on these rows).
$data = @(
    @("s16", "s16_IMG_3179.jpeg", "meltpatch", "1306", "930",  "6", "3", "137", "2"),
    @("s17", "s17_IMG_3175.jpeg", "meltpatch", "2076", "1589", "6", "3", "102", "2"),
    @("s18", "s18_IMG_3174.jpeg", "meltpatch", "2616", "629",  "6", "3", "0",   "2"),
    @("s19", "s19_IMG_3178.jpeg", "meltpatch", "119",  "1543", "6", "3", "70",  "2"),
    @("s20", "s20_IMG_3176.jpeg", "meltpatch", "1219", "900",  "6", "3", "21",  "2")
)

$startRow = 17
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowRange = $ws.Range("A" + $row + ":I" + $row)
    $rowRange.NumberFormat = "@"
    for ($c = 0; $c -lt 9; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $data[$i][$c]
    }
    $rowRange.Style = "Normal"
}
